# Commit: "Added the coulomb potential"
$wb = $excel.ActiveWorkbook

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet4 = $wb.Worksheets.Item("Sheet4")

# --- Fix the "COULOMBCOUPLING" header typo on Sheet2 -> "COULOMB COUPLING" ---
$sheet2.Range("H1").Value = "COULOMB COUPLING"

# --- Sheet3: coulomb coupling column (F2:F7) changes from 1 -> 0.5 ---
$sheet3.Range("F2:F7").Value = 0.5

# --- Sheet4: 2J+1 column (E2:E5 and E7) changes from 2 -> 4 ---
$sheet4.Range("E2").Value = 4
$sheet4.Range("E3").Value = 4
$sheet4.Range("E4").Value = 4
$sheet4.Range("E5").Value = 4
$sheet4.Range("E7").Value = 4

# --- Sheet4: coulomb coupling column (F2:F7) changes from 1 -> -0.5 ---
$sheet4.Range("F2:F7").Value = -0.5

# --- Update each sheet's current selection to match the new cursor positions ---
$sheet2.Activate()
$sheet2.Range("H1").Select()

$sheet3.Activate()
$sheet3.Range("F18").Select()

# Sheet4 becomes the active/visible tab of the workbook
$sheet4.Activate()
$sheet4.Range("F13").Select()
